# Helyesírási javítások, nyelvtani helytelenségek javítva
$d = $word.ActiveDocument

# 1. "felvan tüntetve" -> "fel van tüntetve"
$d.Content.Find.Execute(
    "felvan tüntetve", $false, $false, $false, $false, $false,
    $true, 1, $false, "fel van tüntetve", 2) | Out-Null

# 2. " található ez űrlap, amely bekéri a felhasználó nevét email címét, illetve az üzenetet."
#    -> " található ez az űrlap, amely bekéri a felhasználó nevét, email címét, illetve az üzenetet."
$d.Content.Find.Execute(
    " található ez űrlap, amely bekéri a felhasználó nevét email címét, illetve az üzenetet.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    " található ez az űrlap, amely bekéri a felhasználó nevét, email címét, illetve az üzenetet.",
    2) | Out-Null

# 3. "Asztali számítógép:" -> "Asztali számítógép nézet:"
$d.Content.Find.Execute(
    "Asztali számítógép:", $false, $false, $false, $false, $false,
    $true, 1, $false, "Asztali számítógép nézet:", 2) | Out-Null

# 4. "Tabletes:" -> "Tabletes nézet:"
$d.Content.Find.Execute(
    "Tabletes:", $false, $false, $false, $false, $false,
    $true, 1, $false, "Tabletes nézet:", 2) | Out-Null

# 5. "Telefonos:" -> "Telefonos nézet:"
$d.Content.Find.Execute(
    "Telefonos:", $false, $false, $false, $false, $false,
    $true, 1, $false, "Telefonos nézet:", 2) | Out-Null

# 6. "A legördülő menü-re kattintva" -> "A legördülő menüre kattintva"
$d.Content.Find.Execute(
    "A legördülő menü-re kattintva", $false, $false, $false, $false, $false,
    $true, 1, $false, "A legördülő menüre kattintva", 2) | Out-Null

# 7. "kiszámíttatni az eredményeket" -> "számítani eredményeket"
$d.Content.Find.Execute(
    "kiszámíttatni az eredményeket", $false, $false, $false, $false, $false,
    $true, 1, $false, "számítani eredményeket", 2) | Out-Null

# 8. " akkor a számoló oldal bezáródik és újra" -> " akkor a számoló oldal bezáródik, és újra"
$d.Content.Find.Execute(
    " akkor a számoló oldal bezáródik és újra", $false, $false, $false, $false, $false,
    $true, 1, $false, " akkor a számoló oldal bezáródik, és újra", 2) | Out-Null

# 9. Remove the old _GoBack bookmark (sitting alone in its own empty paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 10. Re-create the _GoBack bookmark collapsed right after the very last
#     run of the document (after "... alakzat választó oldal."). A
#     genuinely empty/collapsed Range at the absolute end of the story
#     cannot host Bookmarks.Add reliably, so a temporary marker character
#     is appended, the bookmark is anchored around it, and the marker is
#     then cleared back out - leaving a zero-length bookmark in place.
$lastPara = $d.Paragraphs.Last
$endPos = $lastPara.Range.End
$insertPoint = $d.Range($endPos - 1, $endPos - 1)
$insertPoint.InsertAfter("Z")

$lastPara2 = $d.Paragraphs.Last
$endPos2 = $lastPara2.Range.End
$marker = $d.Range($endPos2 - 2, $endPos2 - 1)
$d.Bookmarks.Add("_GoBack", $marker)
$marker.Text = ""
